$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Comment text updates on the "Organ" sheet
# ---------------------------------------------------------------------------
$organ = $wb.Worksheets.Item("Organ")

$organ.Range("M1").Comment.Text("The unit of measurement of weight")

$v1Text = "(Required) The string that serves as the definitive identifier for the metadata`nschema version and is readily interpretable by computers for data validation and`nprocessing. Example: 22bc762a-5020-419d-b170-24253ed9e8d9"
$organ.Range("V1").Comment.Text($v1Text)

# ---------------------------------------------------------------------------
# 2. warm_ischemic_time_unit / cold_ischemic_time_unit: drop month/year/day,
#    keep hour + minute only (5 rows -> 2 rows)
# ---------------------------------------------------------------------------
foreach ($name in @("warm_ischemic_time_unit", "cold_ischemic_time_unit")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A1").Value = "hour"
    $ws.Range("B1").Value = "http://purl.obolibrary.org/obo/UO_0000032"
    $ws.Range("A2").Value = "minute"
    $ws.Range("B2").Value = "http://purl.obolibrary.org/obo/UO_0000031"
    $ws.Rows.Item(5).Delete()
    $ws.Rows.Item(4).Delete()
    $ws.Rows.Item(3).Delete()
}

# ---------------------------------------------------------------------------
# 3. weight_unit: drop ng/ug/mg, keep kg + g only (5 rows -> 2 rows)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("weight_unit")
$ws.Range("A1").Value = "kg"
$ws.Range("B1").Value = "http://purl.obolibrary.org/obo/UO_0000009"
$ws.Range("A2").Value = "g"
$ws.Range("B2").Value = "http://purl.obolibrary.org/obo/UO_0000021"
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# 4. height_unit / width_unit / length_unit: drop um/nm, keep mm + cm only
#    (4 rows -> 2 rows)
# ---------------------------------------------------------------------------
foreach ($name in @("height_unit", "width_unit", "length_unit")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A1").Value = "mm"
    $ws.Range("B1").Value = "http://purl.obolibrary.org/obo/UO_0000016"
    $ws.Range("A2").Value = "cm"
    $ws.Range("B2").Value = "http://purl.obolibrary.org/obo/UO_0000015"
    $ws.Rows.Item(4).Delete()
    $ws.Rows.Item(3).Delete()
}

# ---------------------------------------------------------------------------
# 5. volume_unit: reorder to cm^3, um^3, mm^3 and add ml (3 rows -> 4 rows)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("volume_unit")
$ws.Rows.Item(4).Insert()
$ws.Range("A1").Value = "cm^3"
$ws.Range("B1").Value = "http://purl.obolibrary.org/obo/UO_0000097"
$ws.Range("A2").Value = "um^3"
$ws.Range("B2").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000112"
$ws.Range("A3").Value = "mm^3"
$ws.Range("B3").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000136"
$ws.Range("A4").Value = "ml"
$ws.Range("B4").Value = "http://purl.obolibrary.org/obo/UO_0000098"

# ---------------------------------------------------------------------------
# 6. Update data validation ranges on the "Organ" sheet to match the new
#    (shorter / longer) lookup-table extents
# ---------------------------------------------------------------------------
$organ.Range("H2:H1001").Validation.Modify(3, 1, 1, "='warm_ischemic_time_unit'!`$A`$1:`$A`$2")
$organ.Range("J2:J1001").Validation.Modify(3, 1, 1, "='cold_ischemic_time_unit'!`$A`$1:`$A`$2")
$organ.Range("M2:M1001").Validation.Modify(3, 1, 1, "='weight_unit'!`$A`$1:`$A`$2")
$organ.Range("O2:O1001").Validation.Modify(3, 1, 1, "='height_unit'!`$A`$1:`$A`$2")
$organ.Range("Q2:Q1001").Validation.Modify(3, 1, 1, "='width_unit'!`$A`$1:`$A`$2")
$organ.Range("S2:S1001").Validation.Modify(3, 1, 1, "='length_unit'!`$A`$1:`$A`$2")
$organ.Range("U2:U1001").Validation.Modify(3, 1, 1, "='volume_unit'!`$A`$1:`$A`$4")

# ---------------------------------------------------------------------------
# 7. .metadata sheet: bump pav:createdOn timestamp
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item(".metadata")
$meta.Range("C2").Value = "2023-09-08T20:50:05-07:00"
